$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data one column to the right (A:D -> B:E) and one row down (1:11 -> 2:12),
# preserving the per-column widths that were already set (they move along with the columns).
$ws.Columns.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# New column A gets a much wider width to hold the long comparison labels.
$ws.Columns.Item(1).ColumnWidth = 56.166666666666664

# Header row (row 1) for the numeric result columns.
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels in column A describing each comparison.
$ws.Range("A2").Value  = "NbUnique Operands & EffortTo Implement"
$ws.Range("A3").Value  = "NbOperands & EffortTo Implement"
$ws.Range("A4").Value  = "NbUnique Operators & EffortTo Implement"
$ws.Range("A5").Value  = "Program Length & EffortTo Implement"
$ws.Range("A6").Value  = "Difficulty Level & Difficulty Level"
$ws.Range("A7").Value  = "Program Level & Program Level"
$ws.Range("A8").Value  = "EffortTo Implement & NbOperands"
$ws.Range("A9").Value  = "EffortTo Implement & NbUnique Operators"
$ws.Range("A10").Value = "EffortTo Implement & Program Length"
$ws.Range("A11").Value = "EffortTo Implement & EffortTo Implement"
$ws.Range("A12").Value = "TimeTo Implement & TimeTo Implement"
